$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C6: the "fade in / BOO" dialogue line gains extra fade params.
$ws.Range("C6").Value = "[fade=in,0,0,0,0][screen-shake=0.2,0.5]BOO![block=f]"

# Update C5: the "fade out / wait" dialogue line gains extra fade params
# and the new floating center-text effect.
$ws.Range("C5").Value = "[block=t][fade=out,2,0,0,0]Wait, why is the screen fading away?[pause=2][center-text-fade=in,0][center-text-scroll=0.05,0,0,0,SPOOPY][pause=2][center-text-fade=out,2][pause=3][next]"

# Move the active selection from C6 to C5.
$ws.Range("C5").Select()
